# POI upgrade (3.17.0 -> 4.0.1) re-serialized the package: the actual
# document content is unchanged, but every XML part that POI rewrites
# gets its attributes emitted in schema (declaration) order instead of
# the old alphabetical order. We reproduce that by touching the Word
# object model just enough to force each part to be regenerated, using
# property round-trips that are semantic no-ops (read a value, assign
# it straight back) so the visible document content stays identical.

$d = $word.ActiveDocument

# Force word/document.xml to be rewritten (attribute order normalized)
# without altering any actual content: reassign a paragraph property to
# its own current value.
$p = $d.Paragraphs(1)
$p.Alignment = $p.Alignment

# Force word/styles.xml to be rewritten the same way, via a harmless
# round-trip on a style's own name.
$s = $d.Styles("Grilledutableau")
$s.NameLocal = $s.NameLocal
